$d = $word.ActiveDocument

$header = $d.Sections.First.Headers.Item(1)
$header.Range.Find.Execute("2018 Zone 4 and Zone 3/4 Captains", $true, $false, $false, $false, $false,
                            $true, 1, $false, "2018 Zone 4 and Zone 3/4 Team Captains", 2)
